$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# force them to Text format first so the literal string (incl. trailing zeros / formatting) is preserved.
$textCells = @("D5", "D6", "D9", "D10", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D32", "D35", "D37", "D38", "D39", "D40", "D41", "D43", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.565.14"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "3.122.31"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D5").Value = "589.86"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").Value = "145.60"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.119.19"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +15.83%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  +4.87%  "
$ws.Range("D14").Value = "36.10"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "3.649.29"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "7.17"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "63.515.07"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "3.122.96"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "464.10"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "14.29"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").Value = "0.735"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").Value = "82.40"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  +9.02%  "
$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "27.13"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").Value = "0.0₃0851"
$ws.Range("E34").Value = "  +5.89%  "
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  +8.87%  "
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").Value = "3.35"
$ws.Range("E37").Value = "  +11.89%  "
$ws.Range("D38").Value = "6.10"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "50.80"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "446.48"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").Value = "8.75"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.914.38"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0370"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "35.45"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "125.31"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "24.69"
$ws.Range("E51").Value = "  +2.90%  "
